$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete row 2 (the "2024-09-04" row). This shifts all subsequent rows up by one,
#    which matches the diff's row-by-row value shift pattern (J2 takes the old J3
#    value, J3 takes the old J4 value, etc.).
$ws.Rows(2).Delete()

# 2. Determine the new last row (now row 29, holding "2024-10-02").
$lastRow = $ws.UsedRange.Rows.Item($ws.UsedRange.Rows.Count).Row
$newRow = $lastRow + 1

# 3. Copy the date cell's format (style) into the new row so it keeps the same
#    bordered/bold/centered style as every other date cell in column A.
$ws.Cells.Item($lastRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)  # xlPasteFormats

# 4. Append a brand-new row for "2024-10-03", copying the numeric columns (B:J)
#    from the current last row (identical values, as shown by the diff having no
#    further per-cell changes beyond the row shift).
$ws.Cells.Item($lastRow, 2).Resize(1, 9).Copy()
$ws.Cells.Item($newRow, 2).PasteSpecial(-4163)  # xlPasteValues

# 5. Set the date label for the new row as text (not an auto-converted date
#    serial), without disturbing styles. Build the text via a helper formula
#    cell elsewhere, then paste-special (values only) into place so the target
#    cell keeps its style (s="1") and the shared-string table simply gains the
#    new "2024-10-03" entry.
$helperRow = $newRow + 50
$helper = $ws.Cells.Item($helperRow, 1)
$helper.Formula = "=""2024-10-03"""
$helper.Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163)  # xlPasteValues
$helper.Clear()

$excel.CutCopyMode = 0
